# Updated symbol list (cryptos.xlsx) - refresh Price (D) / Volume(1h) (E)
# columns for the scraped coin rows. Values are written with a leading
# apostrophe so Excel stores them as literal text (matching the original
# inlineStr cells) instead of re-interpreting numeric-looking strings
# (e.g. "310.76") or percentages (e.g. "2.12%") as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'310.76"
$ws.Range("E2").Value = "'2.12%"

$ws.Range("D3").Value = "'38.88"
$ws.Range("E3").Value = "'8.59%"

$ws.Range("D4").Value = "'5.117"
$ws.Range("E4").Value = "'1.61%"

$ws.Range("D5").Value = "'0.08175"
$ws.Range("E5").Value = "'2.73%"

$ws.Range("D6").Value = "'2.012"
$ws.Range("E6").Value = "'7.49%"

$ws.Range("D7").Value = "'7.931"
$ws.Range("E7").Value = "'2.07%"

$ws.Range("D8").Value = "'0.9334"
$ws.Range("E8").Value = "'1.43%"

$ws.Range("D9").Value = "'0.1408"
$ws.Range("E9").Value = "'5.15%"

$ws.Range("E10").Value = "'3.47%"

$ws.Range("D11").Value = "'0.09176"
$ws.Range("E11").Value = "'0.57%"

$ws.Range("D12").Value = "'0.03455"
$ws.Range("E12").Value = "'0.53%"

$ws.Range("E13").Value = "'0.00%"

$ws.Range("D14").Value = "'0.001419"
$ws.Range("E14").Value = "'1.38%"

$ws.Range("D15").Value = "'0.005889"
$ws.Range("E15").Value = "'-2.78%"

$ws.Range("E16").Value = "'-4.63%"

$ws.Range("D17").Value = "'4.195"
$ws.Range("E17").Value = "'1.84%"

$ws.Range("D18").Value = "'3.439"
$ws.Range("E18").Value = "'1.33%"

$ws.Range("E19").Value = "'0.22%"

$ws.Range("D20").Value = "'0.1314"
$ws.Range("E20").Value = "'0.34%"

$ws.Range("D21").Value = "'4.835"
$ws.Range("E21").Value = "'-6.16%"

$ws.Range("D22").Value = "'0.2467"
$ws.Range("E22").Value = "'5.05%"

$ws.Range("D23").Value = "'0.04465"
$ws.Range("E23").Value = "'1.31%"

$ws.Range("D24").Value = "'0.001236"
$ws.Range("E24").Value = "'0.40%"

$ws.Range("E25").Value = "'-9.81%"

$ws.Range("E27").Value = "'4.21%"

$ws.Range("D39").Value = "'0.02142"
$ws.Range("E39").Value = "'10.19%"

$ws.Range("D40").Value = "'0.05192"
$ws.Range("E40").Value = "'-3.55%"

$ws.Range("D41").Value = "'0.007476"
$ws.Range("E41").Value = "'-1.66%"

$ws.Range("D42").Value = "'0.009973"
$ws.Range("E42").Value = "'-1.44%"

$ws.Range("D43").Value = "'0.1369"
$ws.Range("E43").Value = "'1.25%"

$ws.Range("D44").Value = "'0.002133"
$ws.Range("E44").Value = "'-1.19%"

$ws.Range("E45").Value = "'-0.38%"

$ws.Range("D46").Value = "'0.00006335"
$ws.Range("E46").Value = "'3.55%"

$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'0.19%"

$ws.Range("E48").Value = "'-0.63%"

$ws.Range("D49").Value = "'0.001602"
$ws.Range("E49").Value = "'-3.31%"

$ws.Range("D50").Value = "'0.00002103"
$ws.Range("E50").Value = "'0.19%"

$ws.Range("D51").Value = "'0.0002003"
$ws.Range("E51").Value = "'0.19%"
